# Update "想去人数" (column F, "people interested") counts on all 4 sheets
# to match the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value  = 1479
$ws.Cells.Item(5, 6).Value  = 7386
$ws.Cells.Item(6, 6).Value  = 73
$ws.Cells.Item(7, 6).Value  = 4755
$ws.Cells.Item(8, 6).Value  = 7001
$ws.Cells.Item(9, 6).Value  = 12
$ws.Cells.Item(10, 6).Value = 268
$ws.Cells.Item(11, 6).Value = 1481
$ws.Cells.Item(12, 6).Value = 853
$ws.Cells.Item(13, 6).Value = 165
$ws.Cells.Item(14, 6).Value = 49
$ws.Cells.Item(17, 6).Value = 156
$ws.Cells.Item(19, 6).Value = 221
$ws.Cells.Item(20, 6).Value = 31
$ws.Cells.Item(21, 6).Value = 1152
$ws.Cells.Item(24, 6).Value = 45
$ws.Cells.Item(25, 6).Value = 1219
$ws.Cells.Item(26, 6).Value = 42
$ws.Cells.Item(27, 6).Value = 141
$ws.Cells.Item(29, 6).Value = 41
$ws.Cells.Item(30, 6).Value = 166
$ws.Cells.Item(32, 6).Value = 34
$ws.Cells.Item(33, 6).Value = 87
$ws.Cells.Item(34, 6).Value = 31
$ws.Cells.Item(37, 6).Value = 66
$ws.Cells.Item(39, 6).Value = 366
$ws.Cells.Item(40, 6).Value = 1196
$ws.Cells.Item(41, 6).Value = 574
$ws.Cells.Item(42, 6).Value = 139
$ws.Cells.Item(43, 6).Value = 19

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value  = 16
$ws.Cells.Item(23, 6).Value = 140
$ws.Cells.Item(28, 6).Value = 24
$ws.Cells.Item(31, 6).Value = 847
$ws.Cells.Item(36, 6).Value = 110
$ws.Cells.Item(43, 6).Value = 73

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 726
$ws.Cells.Item(6, 6).Value = 666
$ws.Cells.Item(8, 6).Value = 1555
$ws.Cells.Item(9, 6).Value = 2446

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value  = 726
$ws.Cells.Item(3, 6).Value  = 1479
$ws.Cells.Item(6, 6).Value  = 666
$ws.Cells.Item(7, 6).Value  = 666
$ws.Cells.Item(8, 6).Value  = 7386
$ws.Cells.Item(9, 6).Value  = 73
$ws.Cells.Item(10, 6).Value = 4755
$ws.Cells.Item(12, 6).Value = 7001
$ws.Cells.Item(13, 6).Value = 268
$ws.Cells.Item(14, 6).Value = 1481
$ws.Cells.Item(16, 6).Value = 853
$ws.Cells.Item(17, 6).Value = 165
$ws.Cells.Item(18, 6).Value = 1555
$ws.Cells.Item(19, 6).Value = 2446
$ws.Cells.Item(21, 6).Value = 49
$ws.Cells.Item(23, 6).Value = 140
$ws.Cells.Item(24, 6).Value = 156
$ws.Cells.Item(25, 6).Value = 221
$ws.Cells.Item(26, 6).Value = 1152
$ws.Cells.Item(29, 6).Value = 1219
$ws.Cells.Item(30, 6).Value = 141
$ws.Cells.Item(31, 6).Value = 166
$ws.Cells.Item(32, 6).Value = 24
$ws.Cells.Item(33, 6).Value = 847
$ws.Cells.Item(34, 6).Value = 34
$ws.Cells.Item(35, 6).Value = 87
$ws.Cells.Item(39, 6).Value = 66
$ws.Cells.Item(41, 6).Value = 110
$ws.Cells.Item(42, 6).Value = 366
$ws.Cells.Item(43, 6).Value = 574
$ws.Cells.Item(47, 6).Value = 139
$ws.Cells.Item(48, 6).Value = 73
$ws.Cells.Item(49, 6).Value = 19
